# Insert a new weekly price record as row 43 in the "Naranja" sheet.
# This pushes the existing rows 43:131 down to 44:132 (dimension becomes
# A1:T132) and populates the newly inserted row 43 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 43 and below down by one row.
$ws.Rows("43:43").Insert()

# Populate the new row 43 with the new price observation.
$ws.Cells.Item(43, 1).Value  = 1
$ws.Cells.Item(43, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(43, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(43, 4).Value  = 44994
$ws.Cells.Item(43, 5).Value  = 15
$ws.Cells.Item(43, 6).Value  = "Fruta"
$ws.Cells.Item(43, 7).Value  = 100102
$ws.Cells.Item(43, 8).Value  = "Cítricos"
$ws.Cells.Item(43, 9).Value  = 100102005
$ws.Cells.Item(43, 10).Value = "Naranja"
$ws.Cells.Item(43, 11).Value = "Valencia"
$ws.Cells.Item(43, 12).Value = "Tercera"
$ws.Cells.Item(43, 13).Value = 270
$ws.Cells.Item(43, 14).Value = 900
$ws.Cells.Item(43, 15).Value = 1000
$ws.Cells.Item(43, 16).Value = 950
$ws.Cells.Item(43, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(43, 18).Value = "Región Metropolitana"
$ws.Cells.Item(43, 19).Value = 950
$ws.Cells.Item(43, 20).Value = 1
